$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.270.95"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.073.44"
$ws.Range("E3").Value = "  +4.07%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.54"
$ws.Range("E5").Value = "  -2.25%  "
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.19"
$ws.Range("E8").Value = "  +4.18%  "
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.03"
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("E12").Value = "  +3.18%  "
$ws.Range("D13").Value = "2.376.21"
$ws.Range("E13").Value = "  +3.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.52"
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.05"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.23"
$ws.Range("E17").Value = "  +3.60%  "
$ws.Range("D18").Value = "2.098.43"
$ws.Range("E18").Value = "  +5.02%  "
$ws.Range("D19").Value = "37.384.63"
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.96"
$ws.Range("E20").Value = "  +19.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "68.37"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "0.0₃0811"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "223.89"
$ws.Range("E23").Value = "  -2.07%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.41"
$ws.Range("E25").Value = "  +2.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.43"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.61"
$ws.Range("E27").Value = "  +0.78%  "
$ws.Range("E28").Value = "  +2.35%  "
$ws.Range("E29").Value = "  +5.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.24"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  +7.73%  "
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.44"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0623"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.53"
$ws.Range("E35").Value = "  +8.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.39"
$ws.Range("E36").Value = "  +4.09%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +13.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.32"
$ws.Range("E39").Value = "  +0.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.76"
$ws.Range("E40").Value = "  -1.30%  "
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("B42").Value = "FTXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.46"
$ws.Range("E42").Value = "  +22.80%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0962"
$ws.Range("E43").Value = "  +8.60%  "
$ws.Range("D44").Value = "1.471.05"
$ws.Range("E44").Value = "  +2.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "94.88"
$ws.Range("E45").Value = "  +7.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0209"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "16.12"
$ws.Range("E47").Value = "  +5.12%  "
$ws.Range("E48").Value = "  +0.83%  "
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.16"
$ws.Range("E50").Value = "  +8.04%  "
$ws.Range("E51").Value = "  +1.79%  "
